# New API Query - 2023 Included
# API query to UN performed 11/26/2023.
# Query modified to include 2023 data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "short-url" id used for this query changed for every data row.
$ws.Range("B2:B5").Value = "hYbS1u"

# Existing "0" placeholders for the "hst" (host) column become "-" for the
# first three data rows, and pick up the left-aligned style used for
# non-numeric placeholder text.
$ws.Range("V2:V4").Value = "-"
$ws.Range("V2:V4").HorizontalAlignment = -4131

# Row 5's "oip" value of "null" becomes "-" as well (already left-aligned).
$ws.Range("U5").Value = "-"

# Append a new data row (row 6) for year 2023.
$ws.Range("A6").Value = "1"
$ws.Range("A6").HorizontalAlignment = -4152

$ws.Range("B6").Value = "hYbS1u"
$ws.Range("B6").HorizontalAlignment = -4131

$ws.Range("C6").Value = "1"
$ws.Range("C6").HorizontalAlignment = -4152

$ws.Range("D6").Value = "5"
$ws.Range("D6").HorizontalAlignment = -4152

$ws.Range("E6").Value = "2023"
$ws.Range("E6").HorizontalAlignment = -4152

$ws.Range("F6").Value = "207"
$ws.Range("F6").HorizontalAlignment = -4152

$ws.Range("G6").Value = "Venezuela (Bolivarian Republic of)"
$ws.Range("G6").HorizontalAlignment = -4131

$ws.Range("H6").Value = "VEN"
$ws.Range("H6").HorizontalAlignment = -4131

$ws.Range("I6").Value = "VEN"
$ws.Range("I6").HorizontalAlignment = -4131

$ws.Range("J6").Value = "206"
$ws.Range("J6").HorizontalAlignment = -4152

$ws.Range("K6").Value = "Saint Vincent and the Grenadines"
$ws.Range("K6").HorizontalAlignment = -4131

$ws.Range("L6").Value = "VCT"
$ws.Range("L6").HorizontalAlignment = -4131

$ws.Range("M6").Value = "VCT"
$ws.Range("M6").HorizontalAlignment = -4131

$ws.Range("N6").Value = "0"
$ws.Range("N6").HorizontalAlignment = -4152

$ws.Range("O6").Value = "0"
$ws.Range("O6").HorizontalAlignment = -4152

$ws.Range("P6").Value = "0"
$ws.Range("P6").HorizontalAlignment = -4152

$ws.Range("Q6").Value = "0"
$ws.Range("Q6").HorizontalAlignment = -4152

$ws.Range("R6").Value = "0"
$ws.Range("R6").HorizontalAlignment = -4152

$ws.Range("S6").Value = "0"
$ws.Range("S6").HorizontalAlignment = -4152

$ws.Range("T6").Value = "7"
$ws.Range("T6").HorizontalAlignment = -4152

$ws.Range("U6").Value = "-"
$ws.Range("U6").HorizontalAlignment = -4131

$ws.Range("V6").Value = "0"
$ws.Range("V6").HorizontalAlignment = -4152
